# Generate Report for Handback
# Updates the timestamp values recorded in the handback-status report.
#
# Mapping of changed timestamps (old -> new):
#   2016-08-25 17:07:43 -> 2016-08-25 17:08:50   (Overview sheet G2, "Latest HO Xliff Generate Date";
#                                                  shared with de-de sheet H2, "Correspond Handoff Datetime")
#   2016-08-25 17:07:38 -> 2016-08-25 17:08:45   (zh-cn sheet H2, "Correspond Handoff Datetime")
#   2016-08-25 17:08:21 -> 2016-08-25 17:09:27   (zh-cn sheet K2, "Correspond Handback DateTime")
#   2016-08-25 17:08:28 -> 2016-08-25 17:09:35   (de-de sheet K2, "Correspond Handback DateTime")

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" column (G), row 2
# (this text value is also used by de-de!H2, "Correspond Handoff Datetime" -- set both so they
# keep sharing a single, updated string)
$wsOverview.Range("G2").Value = "2016-08-25 17:08:50"
$wsDeDe.Range("H2").Value = "2016-08-25 17:08:50"

# zh-cn: "Correspond Handoff Datetime" column (H) and "Correspond Handback DateTime" column (K), row 2
$wsZhCn.Range("H2").Value = "2016-08-25 17:08:45"
$wsZhCn.Range("K2").Value = "2016-08-25 17:09:27"

# de-de: "Correspond Handback DateTime" column (K), row 2
$wsDeDe.Range("K2").Value = "2016-08-25 17:09:35"
